$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 88 ("Fruta / hortaliza, semanal"),
# pushing all the existing rows 88-185 down to 89-186.
$ws.Rows(88).Insert()

$ws.Cells.Item(88, 1).Value = 8
$ws.Cells.Item(88, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(88, 3).Value = "Coquimbo"
$ws.Cells.Item(88, 4).Value = 45079
$ws.Cells.Item(88, 5).Value = 4
$ws.Cells.Item(88, 6).Value = 100112052
$ws.Cells.Item(88, 7).Value = "Albahaca"
$ws.Cells.Item(88, 8).Value = "Sin especificar"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 900
$ws.Cells.Item(88, 11).Value = 2500
$ws.Cells.Item(88, 12).Value = 3000
$ws.Cells.Item(88, 13).Value = 2750
$ws.Cells.Item(88, 14).Value = "$/paquete"
$ws.Cells.Item(88, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(88, 16).Value = 2750
$ws.Cells.Item(88, 17).Value = 1
$ws.Cells.Item(88, 18).Value = "Hortaliza"
